$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings (e.g. "30.216.50", "95.00")
# that must stay as literal text, matching the source data feed formatting.
# Force text number format before assignment so Excel does not coerce them
# into actual numbers (which would drop significant trailing zeros, etc.).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.216.50"
$ws.Range("E2").Value = "  +5.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.907.48"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.05"
$ws.Range("E5").Value = "  +4.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5165"
$ws.Range("E7").Value = "  +1.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4033"
$ws.Range("E8").Value = "  +3.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08468"
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.76"
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.122"
$ws.Range("E11").Value = "  +1.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.46"
$ws.Range("E12").Value = "  +15.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.419"
$ws.Range("E13").Value = "  +3.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.914.26"
$ws.Range("E14").Value = "  +2.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.371"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "95.00"
$ws.Range("E17").Value = "  +1.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001113"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06711"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.39"
$ws.Range("E20").Value = "  +4.27%  "
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.002"
$ws.Range("E22").Value = "  +1.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.235.99"
$ws.Range("E23").Value = "  +5.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.32"
$ws.Range("E24").Value = "  +2.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.220"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.133.85"
$ws.Range("E26").Value = "  +2.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.84"
$ws.Range("E27").Value = "  +5.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.52"
$ws.Range("E28").Value = "  +3.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.390"
$ws.Range("E29").Value = "  -0.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.75"
$ws.Range("E30").Value = "  +2.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.104"
$ws.Range("E31").Value = "  +5.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1063"
$ws.Range("E32").Value = "  +2.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.996"
$ws.Range("E33").Value = "  +3.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.640"
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02489"
$ws.Range("E35").Value = "  +1.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06573"
$ws.Range("E36").Value = "  +0.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2208"
$ws.Range("E37").Value = "  +2.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.230"
$ws.Range("E38").Value = "  +2.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.170"
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.92"
$ws.Range("E40").Value = "  +6.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.782"
$ws.Range("E41").Value = "  -2.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6514"
$ws.Range("E42").Value = "  +2.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.232"
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6128"
$ws.Range("E44").Value = "  +2.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.29"
$ws.Range("E45").Value = "  +1.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.732"
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.067"
$ws.Range("E47").Value = "  +3.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.244"
$ws.Range("E48").Value = "  +1.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.09"
$ws.Range("E49").Value = "  +2.48%  "
$ws.Range("E50").Value = "  -1.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.38"
$ws.Range("E51").Value = "  +3.99%  "
